# "Generate Report for Handoff" - refresh the localization-status report:
#  - flip Status from "Handed back: in sync with en-US" to "Ready for handoff"
#    on the Overview sheet (zh-cn + de-de columns) and on each language sheet
#  - bump the "Latest Handoff/Generate" timestamp cells to the new run time
#  - the Status column narrows now that the new text is shorter, so re-fit it

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# Overview sheet: E2 = zh-cn status, F2 = de-de status, G2 = generate date
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-10-26 08:11:03"

# zh-cn sheet: C2 = status, H2 = latest handoff datetime
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-10-26 08:10:51"

# de-de sheet: C2 = status, H2 = latest handoff datetime
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-10-26 08:11:03"

# Re-fit the Status columns now that "Ready for handoff" is shorter than the
# previous "Handed back: in sync with en-US" text they held.
$overview.Columns.Item(5).ColumnWidth = 16.25
$overview.Columns.Item(6).ColumnWidth = 16.25
$zhcn.Columns.Item(3).ColumnWidth = 16.25
$dede.Columns.Item(3).ColumnWidth = 16.25
